# Refresh the "cryptos" price list (GitHub Actions data pull).
# Note: several Price-column values are plain decimals (e.g. "1.005"), which
# Excel's Value setter auto-coerces to a number. We force those back to text
# with a leading apostrophe (matching the workbook's original inlineStr
# layout) and then reset Style to "Normal" so no stray NumberFormat/quote
# -prefix formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.189.52"

$ws.Range("D3").Value = "1.608.85"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "'1.004"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'302.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").Value = "'0.3769"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").Value = "'0.3510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.27%  "

$ws.Range("D10").Value = "'0.08049"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("D14").Value = "'6.321"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.49%  "

$ws.Range("D15").Value = "'7.226"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").Value = "'0.00001207"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("D17").Value = "1.591.85"
$ws.Range("E17").Value = "  -2.17%  "

$ws.Range("D18").Value = "'94.05"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'0.06907"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'6.451"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").Value = "'17.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("E23").Value = "  -4.17%  "

$ws.Range("D24").Value = "23.188.34"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("D25").Value = "'2.545"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "

$ws.Range("D26").Value = "'3.044"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.21%  "

$ws.Range("E27").Value = "  -1.65%  "

$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "'5.256"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").Value = "'131.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "

$ws.Range("D31").Value = "1.770.12"
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.057"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.00%  "

$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.34%  "

$ws.Range("D34").Value = "'6.388"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.78%  "

$ws.Range("D35").Value = "'11.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.18%  "

$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("D37").Value = "'0.08658"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("D38").Value = "'0.2444"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").Value = "'0.06863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.72%  "

$ws.Range("D40").Value = "'5.796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.98%  "

$ws.Range("D41").Value = "'0.6805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.69%  "

$ws.Range("D42").Value = "'1.305"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("D43").Value = "'11.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("E44").Value = "  -6.21%  "

$ws.Range("D45").Value = "'1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").Value = "'0.6243"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.20%  "

$ws.Range("D47").Value = "'3.934"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("E48").Value = "  -3.70%  "

$ws.Range("D49").Value = "'0.07844"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("D50").Value = "'127.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.43%  "

$ws.Range("E51").Value = "  -3.37%  "
